# Automatische test-sync: 2025-07-31 21:52:50
# Adds a new row (row 12) to Sheet1 describing Testmail #14 about CE certificates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question = "Heb je de CE-certificaten van dit product?"
$subject = "Testmail #14: " + $question

$answer = "Geachte klant,`nDank u voor uw e-mail. Voor het verkrijgen van de CE-certificaten van het product waar u naar vraagt, verzoeken wij u ons het specifieke productnummer of de productnaam te verstrekken. Met deze informatie kunnen wij u de relevante certificaten verstrekken.`nAls u verdere vragen heeft of meer ondersteuning nodig heeft, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam] E-mailassistent - [Bedrijfsnaam]"

$ws.Range("A12").Value = $subject
$ws.Range("B12").Value = $answer
$ws.Range("C12").Value = $question
$ws.Range("D12").Value = "mailmind.test@zohomail.eu"
$ws.Range("E12").Value = "Productinformatie"
$ws.Range("F12").Value = "2025-07-31 21:52:33"
$ws.Range("G12").Value = "Ja"
$ws.Range("H12").Value = "Nee"
$ws.Range("I12").Value = "Ja"
$ws.Range("J12").Value = "Nee"

# The multi-line answer text triggers Excel's automatic row-height adjustment;
# re-run AutoFit so the row keeps the default (non-custom) height, matching
# the rest of the sheet.
$ws.Rows.Item(12).EntireRow.AutoFit()
